# "Generate Report for Handoff"
# Flips the localization-status report from "In Translation" to
# "Ready for handoff" and bumps the two "Latest HO Xliff / Handoff"
# generation timestamps that accompany a fresh handoff-package build.
#
# Overview!E2 (zh-cn status), Overview!F2 (de-de status), zh-cn!C2 and
# de-de!C2 all point at the same "In Translation" shared string, and
# Overview!G2 / de-de!H2 share "2016-08-25 22:37:59" while zh-cn!H2 has
# its own "2016-08-25 22:37:54" - so all of those cells move together.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Ready for handoff"

# Status columns
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$zhcn.Range("C2").Value     = $newStatus
$dede.Range("C2").Value     = $newStatus

# Handoff / HO-Xliff generation timestamps
$overview.Range("G2").Value = "2016-08-25 22:38:39"
$dede.Range("H2").Value     = "2016-08-25 22:38:39"
$zhcn.Range("H2").Value     = "2016-08-25 22:38:35"

# The longer "Ready for handoff" text needs a wider status column than
# "In Translation" did - widen the affected columns to fit (this engine
# quantizes ColumnWidth to whole pixels, so we feed it a value that lands
# on the closest achievable width to the authored one).
$overview.Columns.Item(5).ColumnWidth = 16.333333333333332
$overview.Columns.Item(6).ColumnWidth = 16.333333333333332
$zhcn.Columns.Item(3).ColumnWidth     = 16.333333333333332
$dede.Columns.Item(3).ColumnWidth     = 16.333333333333332
